$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 234, shifting existing rows (234..285) down to (235..286)
$ws.Rows.Item(234).Insert()

$r = 234
$ws.Cells.Item($r, 1).Value = 10
$ws.Cells.Item($r, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($r, 3).Value = "La Araucanía"
$ws.Cells.Item($r, 4).Value2 = 44642
$ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 5).Value = 9
$ws.Cells.Item($r, 6).Value = "Fruta"
$ws.Cells.Item($r, 7).Value = 100103
$ws.Cells.Item($r, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item($r, 9).Value = 100103004
$ws.Cells.Item($r, 10).Value = "Durazno"
$ws.Cells.Item($r, 11).Value = "September Sweet"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 90
$ws.Cells.Item($r, 14).Value = 18000
$ws.Cells.Item($r, 15).Value = 18000
$ws.Cells.Item($r, 16).Value = 18000
$ws.Cells.Item($r, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item($r, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($r, 19).Value = 1000
$ws.Cells.Item($r, 20).Value = 18
